$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$cell = $table.Cell(2, 2)
$rng = $cell.Range

$rng.Find.Execute("api/doctor/filter", $true, $false, $false, $false, $false,
                   $true, 1, $false, "api/admin/login", 2) | Out-Null
Write-Host "After find/replace: Start=$($rng.Start) End=$($rng.End)"

# Do a brand NEW find, scanning from the whole doc content, searching for
# the literal replaced text, to get an independent/fresh Range untouched by
# the previous Find's internal state.
$rng2 = $d.Content
$rng2.Find.Execute("api/admin/login", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0) | Out-Null
Write-Host "rng2 after fresh find: Start=$($rng2.Start) End=$($rng2.End)"
$rng2.Collapse(0)
Write-Host "rng2 after collapse: Start=$($rng2.Start) End=$($rng2.End)"

$d.Bookmarks.Add("TB3", $rng2) | Out-Null
$bm = $d.Bookmarks.Item("TB3")
Write-Host "TB3 Start=$($bm.Start) End=$($bm.End)"
